# Corrected Main Class Diagram for DG
#
# 1) Fix the typo "CrearCommand" -> "ClearCommand" in the class diagram
#    rectangle on slide 3.
# 2) Bump the cached date-field text from 11/11/2018 to 11/12/2018
#    everywhere it appears: the Slide Master, every Custom Layout, and
#    the Notes Master.

$p = $ppt.ActivePresentation

# --- 1) Fix "CrearCommand" -> "ClearCommand" on slide 3 ---------------
$slide = $p.Slides.Item(3)
$shape = $slide.Shapes.Item(31)
if ($shape.TextFrame.TextRange.Text -eq "CrearCommand") {
    $shape.TextFrame.TextRange.Text = "ClearCommand"
}

# --- 2) Update cached date field text from 11/11/2018 to 11/12/2018 ---
$oldDate = "11/11/2018"
$newDate = "11/12/2018"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DateShape $p.SlideMaster.Shapes

# Every Custom Layout hanging off the Slide Master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShape $layouts.Item($j).Shapes
}

# Notes Master
if ($p.HasNotesMaster) {
    Update-DateShape $p.NotesMaster.Shapes
}
